# PARGT_Manual_Windows.docx edit:
#   "To be submitted. PARGT: A Standalone Software Tool for Predicting
#    Antimicrobial Resistance in Bacteria." ->
#   "To be submitted. PARGT: A Software Tool for Predicting Antimicrobial
#    Resistance in Bacteria."
# i.e. drop the stray word "Standalone " from the citation entry, and
# (as Word automatically does) the "_GoBack" last-edit bookmark follows the
# edit from its old spot (by "predicted_resistance_sequences.fasta") to the
# new edit location right after "PARGT: A ".

$d = $word.ActiveDocument

# Locate "PARGT: A Standalone Software Tool" and capture its character span.
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "PARGT: A Standalone Software Tool", $false, $false, $false,
    $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $findRange.Start
    $end = $findRange.End

    # Position right after "PARGT: A " (before "Standalone").
    $bookmarkPos = $start + 9

    # Remove "Standalone " from that exact span only (keeps the edit
    # scoped so unrelated runs elsewhere aren't touched).
    $scopedRange = $d.Range($start, $end)
    $scopedRange.Find.Execute(
        "Standalone ", $false, $false, $false, $false, $false,
        $true, 1, $false, "", 2)

    # Move the "_GoBack" bookmark to the new edit location (collapsed,
    # zero-length range). Re-adding a bookmark with an existing name
    # relocates it rather than creating a duplicate, which also removes
    # it from its old spot automatically.
    $bmRange = $d.Range($bookmarkPos, $bookmarkPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
